# Update TC_ID Excel SCD0017 until SCD0025 and Update TC_ID Solution SCD0006 until SCD0025
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from SCD0325 to SCD0024
$ws.Name = "SCD0024"

# Update the TC_ID value in B2 from "DGS-340" to "SCD0024-004"
$ws.Range("B2").Value = "SCD0024-004"

# Adjust alignment of the row-2 header cells (A2, B2, C2, I2) to horizontal=left
$ws.Range("A2").HorizontalAlignment = -4131   # xlLeft
$ws.Range("B2").HorizontalAlignment = -4131   # xlLeft
$ws.Range("B2").VerticalAlignment = -4108     # xlCenter
$ws.Range("C2").HorizontalAlignment = -4131   # xlLeft
$ws.Range("I2").HorizontalAlignment = -4131   # xlLeft

# Widen column B to fit the new longer TC_ID text
$ws.Columns("B").ColumnWidth = 11.5

# Move the active selection to B3
$null = $ws.Range("B3").Select()
